$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.147.21"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "2.051.87"
$ws.Range("E3").Value = "  -1.72%  "
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'248.36"
$ws.Range("E5").Value = "  -2.62%  "
$ws.Range("D7").Value = "'58.31"
$ws.Range("E7").Value = "  -5.80%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -2.94%  "
$ws.Range("E10").Value = "  -3.26%  "
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("D13").Value = "2.349.84"
$ws.Range("E13").Value = "  -1.63%  "
$ws.Range("D14").Value = "'0.839"
$ws.Range("E14").Value = "  +1.24%  "
$ws.Range("D15").Value = "'5.72"
$ws.Range("E15").Value = "  +2.70%  "
$ws.Range("D16").Value = "2.052.55"
$ws.Range("E16").Value = "  -1.67%  "
$ws.Range("D17").Value = "'18.06"
$ws.Range("E17").Value = "  +14.76%  "
$ws.Range("D18").Value = "37.208.55"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").Value = "'74.87"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").Value = "0.0₃0897"
$ws.Range("E20").Value = "  -4.21%  "
$ws.Range("E21").Value = "  -2.87%  "
$ws.Range("D22").Value = "'237.21"
$ws.Range("E22").Value = "  -1.71%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  +1.52%  "
$ws.Range("D25").Value = "'2.18"
$ws.Range("E25").Value = "  -6.97%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'9.47"
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'169.44"
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("D28").Value = "'20.07"
$ws.Range("E28").Value = "  -2.14%  "
$ws.Range("E29").Value = "  -2.02%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "'4.80"
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "'1.11"
$ws.Range("E31").Value = "  -2.25%  "
$ws.Range("E32").Value = "  -3.94%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").Value = "'0.0894"
$ws.Range("E34").Value = "  -2.97%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").Value = "  -3.23%  "
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").Value = "'3.25"
$ws.Range("E38").Value = "  +15.10%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'1.34"
$ws.Range("E39").Value = "  -2.81%  "
$ws.Range("D40").Value = "'5.21"
$ws.Range("E40").Value = "  +15.03%  "
$ws.Range("E41").Value = "  -16.90%  "
$ws.Range("E42").Value = "  -3.02%  "
$ws.Range("D43").Value = "'17.17"
$ws.Range("E43").Value = "  -6.64%  "
$ws.Range("E44").Value = "  -3.20%  "
$ws.Range("D45").Value = "'95.90"
$ws.Range("E45").Value = "  -4.37%  "
$ws.Range("E46").Value = "  -2.81%  "
$ws.Range("D47").Value = "1.277.18"
$ws.Range("E47").Value = "  -3.07%  "
$ws.Range("D48").Value = "'2.85"
$ws.Range("E48").Value = "  -4.55%  "
$ws.Range("E49").Value = "  -2.57%  "
$ws.Range("D50").Value = "2.236.40"
$ws.Range("E50").Value = "  -1.75%  "
$ws.Range("B51").Value = "FTXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D51").Value = "'3.55"
$ws.Range("E51").Value = "  -21.61%  "
